# ---------------------------------------------------------------------------
# post-storm-content.xlsx -> "post storm" content refresh (Hurricane Josh)
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "content" sheet values (column B) with the new copy.
#    Column A keys are unchanged; only which shared string each B cell
#    points at changes (several old "placeholder" strings are retired and
#    replaced by the real Hurricane-Josh strings below).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("content")

$ws.Range("B3").Value  = "Hurricane Josh Shelter Locator"
$ws.Range("B4").Value  = "If you need emergency shelter as a result of Hurricane Josh, use the map to find the nearest facility or call 311 for info"
$ws.Range("B5").Value  = "Shelter Map"
$ws.Range("B6").Value  = "Shelters"
$ws.Range("B7").Value  = "Shelters"
$ws.Range("B8").Value  = "If your home is unsafe, you may go to one of the shelters listed below. If you cannot get there on your own please call 311. CHECK BACK DAILY AS THE LIST OF SHELTERS MAY CHANGE"
$ws.Range("B9").Value  = "Areas impacted by Hurricane Josh may still be unsafe. If you cannot remain in your home, use this application or call 311 to located a city-operated shelter."
$ws.Range("B10").Value = "shelter"
$ws.Range("B11").Value = "You are not located in an area impaced by Hurricane Josh storm surge"
$ws.Range("B12").Value = "Zone Finder cannot determine your address.<br>Try alternative address or examine map and click on your location."
$ws.Range("B13").Value = "Zone Finder cannot determine your address.<br>Try alternative address."
$ws.Range("B14").Value = "Hurricane Josh Shelter Locator"
$ws.Range("B15").Value = "You are in an area that may have been impacted by Hurricane Josh"
$ws.Range("B16").ClearContents()

$ws.Range("B16").Select()

# ---------------------------------------------------------------------------
# 2) Add the hidden Esri bookkeeping sheet at the end of the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$esri = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$esri.Name = "ESRI_MAPINFO_SHEET"

$shp = $esri.Shapes.AddTextbox(1, 0, 0, 503, 130)
$shp.Name = "EsriDoNotEdit"
$chars = $shp.TextFrame.Characters()
$chars.Text = "DO NOT EDIT `r For Esri use only"
$chars.Font.Size = 50
$chars.Font.Bold = $true
$chars.Font.Name = "Verdana"

$esri.Visible = [Microsoft.Office.Interop.Excel.XlSheetVisibility]::xlSheetVeryHidden

# Re-select the content sheet as the active tab when done.
$ws.Activate()
